$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 220 (shifts existing rows 220:233 down to 221:234),
# mirroring Excel's Rows("220:220").Insert behaviour.
$ws.Rows("220:220").Insert()

# Populate the newly inserted row 220 with the new weekly price record.
# (A new row is blank after Insert, so every column must be (re)written -
# most values repeat the market/category metadata that was already present
# on the old row 220, with D, J, K, L, M, P holding the new figures.)
$ws.Cells.Item(220, 1).Value = 7                                    # A - Mercado ID
$ws.Cells.Item(220, 2).Value = "Terminal Hortofrutícola Agro Chillán" # B - Mercado
$ws.Cells.Item(220, 3).Value = "Ñuble"                               # C - Región
$ws.Cells.Item(220, 4).Value = 45267                                 # D - Fecha (2023-12-07)
$ws.Cells.Item(220, 5).Value = 16                                    # E - Codreg
$ws.Cells.Item(220, 6).Value = 100112037                             # F - Categoría ID
$ws.Cells.Item(220, 7).Value = "Cebollín"                            # G - Categoría
$ws.Cells.Item(220, 8).Value = "Sin especificar"                    # H - Variedad
$ws.Cells.Item(220, 9).Value = "Primera"                             # I - Calidad
$ws.Cells.Item(220, 10).Value = 120                                  # J - Volumen
$ws.Cells.Item(220, 11).Value = 5000                                 # K - Precio minimo
$ws.Cells.Item(220, 12).Value = 5000                                 # L - Precio maximo
$ws.Cells.Item(220, 13).Value = 5000                                 # M - Precio promedio ponderado
$ws.Cells.Item(220, 14).Value = "$/paquete 36 unidades"              # N - Unidad de comercializacion
$ws.Cells.Item(220, 15).Value = "Provincia de Diguillín"             # O - Origen
$ws.Cells.Item(220, 16).Value = 139                                  # P - Precio $/Kg
$ws.Cells.Item(220, 17).Value = 36                                   # Q - Kg o Unidades
$ws.Cells.Item(220, 18).Value = "Hortaliza"                          # R - Clasificacion
